# Edit script: update "烧碱(折100%)" sheet
# 1. Remove columns F (产销率) and G (销售量) entirely - they duplicated data
#    that already exists in columns B (产销率_累计值) and E (销售量_累计值).
# 2. Within each year, the "B" and "C" quarter rows had been listed out of
#    chronological order; swap them back into order (B before C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($row1, $row2)
    $range1 = $ws.Range("A" + $row1 + ":E" + $row1)
    $range2 = $ws.Range("A" + $row2 + ":E" + $row2)
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Swap the mis-ordered quarter rows for each year (rows are 1-indexed,
# row 1 is the header).
Swap-Rows 3 4
Swap-Rows 7 8
Swap-Rows 11 12
Swap-Rows 15 16

# Remove the redundant F (产销率) and G (销售量) columns.
$ws.Range("F1:G17").Delete()
